$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Crypto price/volume refresh values (from upstream data source)
# D-column price values are written as text to avoid Excel auto-parsing
# strings like "601.78" as numbers (values with multiple dots such as
# "68.751.98" are thousands-grouped prices already stored as text).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.751.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.70%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.733.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.96%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.64%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.22"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.48%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.733.29"
$ws.Range("D7").Style = "Normal"

$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.536"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.83%  "

$ws.Range("E10").Value = "  +2.40%  "

$ws.Range("E11").Value = "  +2.61%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.458"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.08%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.96"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.00%  "

$ws.Range("E14").Value = "  -0.36%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.358.93"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.90%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.734.97"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.99%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.747.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.61%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.23"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.43%  "

$ws.Range("E19").Value = "  +0.43%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.20%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "496.95"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.10%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +10.49%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.722"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.83%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.15%  "

$ws.Range("E25").Value = "  -2.82%  "

$ws.Range("E26").Value = "  -7.20%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.36"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.47%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.02%  "

$ws.Range("E29").Value = "  -0.02%  "

$ws.Range("E30").Value = "  -0.80%  "

$ws.Range("E31").Value = "  -0.20%  "

$ws.Range("E32").Value = "  +3.15%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.57"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.16%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.879.83"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.78%  "

$ws.Range("E35").Value = "  -0.60%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.666.66"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.09%  "

$ws.Range("E37").Value = "  +0.05%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.01"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.10%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.80"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.12%  "

$ws.Range("E40").Value = "  -2.29%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.324"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.75%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "432.69"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.97%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "49.09"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.17%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.97"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.51%  "

$ws.Range("E45").Value = "  -1.18%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.40"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.78%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.56"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.62%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.74"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.48%  "

$ws.Range("E50").Value = "  +0.98%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.741.66"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.65%  "
